$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits right after "DUBAI, UAE."
#    It needs to move down into the "NET WEIGHT" paragraph, so remove
#    it from its old spot first.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Find the misspelled "MOSISTURE" inside the NET WEIGHT paragraph.
# ------------------------------------------------------------------
$f = $d.Content.Duplicate
[void]$f.Find.Execute("MOSISTURE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$foundStart = $f.Start

# ------------------------------------------------------------------
# 3. Remove the extra "S" right after "MO" so the word reads
#    "MOISTURE" instead of "MOSISTURE".
# ------------------------------------------------------------------
$badS = $d.Range($foundStart + 2, $foundStart + 3)
$badS.Text = ""

# ------------------------------------------------------------------
# 4. Keep "MO" and "I" as distinct runs (matching the authored
#    edit) by marking the boundary with a temporary bookmark, then
#    drop the _GoBack bookmark right after "I" (before "STURE").
# ------------------------------------------------------------------
$splitPoint = $d.Range($foundStart + 2, $foundStart + 2)
$d.Bookmarks.Add("ZZZ_TempSplit", $splitPoint)

$bmPos = $foundStart + 3
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Bookmarks.Item("ZZZ_TempSplit").Delete()

Write-Host "Paragraph now reads: $($d.Paragraphs.Item(14).Range.Text)"
